# Updates res_bus/vm_pu.xlsx values for the "case with 380 kV" run.
# Replaces the per-bus voltage magnitude (p.u.) results in columns C:F and J:N
# for data rows 2-25 with the recomputed values from the new case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (bus index 0)
$ws.Range("C2").Value = 1.040175243839556
$ws.Range("D2").Value = 1.0428196360626
$ws.Range("E2").Value = 1.047274954924093
$ws.Range("F2").Value = 1.055595767166296
$ws.Range("J2").Value = 1.045263271193734
$ws.Range("K2").Value = 1.045595246361179
$ws.Range("L2").Value = 1.050038048749574
$ws.Range("M2").Value = 1.058335791655318
$ws.Range("N2").Value = 1.018349091485198

# Row 3 (bus index 1)
$ws.Range("C3").Value = 1.042716868103573
$ws.Range("D3").Value = 1.04527617680278
$ws.Range("E3").Value = 1.049600976006425
$ws.Range("F3").Value = 1.058140979816293
$ws.Range("J3").Value = 1.047442465382062
$ws.Range("K3").Value = 1.047858394805339
$ws.Range("L3").Value = 1.052171958826295
$ws.Range("M3").Value = 1.060690067838608
$ws.Range("N3").Value = 1.019131940075214

# Row 4 (bus index 2)
$ws.Range("C4").Value = 1.044352909264869
$ws.Range("D4").Value = 1.046857622668442
$ws.Range("E4").Value = 1.051098044896565
$ws.Range("F4").Value = 1.05977984141871
$ws.Range("J4").Value = 1.048844114537247
$ws.Range("K4").Value = 1.049314436128667
$ws.Range("L4").Value = 1.053544418375485
$ws.Range("M4").Value = 1.062205120318036
$ws.Range("N4").Value = 1.019634041792785

# Row 5 (bus index 3)
$ws.Range("C5").Value = 1.045038707707896
$ws.Range("D5").Value = 1.04752057609493
$ws.Range("E5").Value = 1.051725542078274
$ws.Range("F5").Value = 1.060466943123821
$ws.Range("J5").Value = 1.049431397062961
$ws.Range("K5").Value = 1.049924602763515
$ws.Range("L5").Value = 1.054119455211515
$ws.Range("M5").Value = 1.06284010745029
$ws.Range("N5").Value = 1.019844076417721

# Row 6 (bus index 4)
$ws.Range("C6").Value = 1.045153741009737
$ws.Range("D6").Value = 1.047631779792947
$ws.Range("E6").Value = 1.051830793363836
$ws.Range("F6").Value = 1.06058220209053
$ws.Range("J6").Value = 1.049529890216212
$ws.Range("K6").Value = 1.050026939332132
$ws.Range("L6").Value = 1.054215893725981
$ws.Range("M6").Value = 1.062946612243975
$ws.Range("N6").Value = 1.019879281142692

# Row 7 (bus index 5)
$ws.Range("C7").Value = 1.044362080694963
$ws.Range("D7").Value = 1.046866488426631
$ws.Range("E7").Value = 1.051106436828204
$ws.Range("F7").Value = 1.059789029803533
$ws.Range("J7").Value = 1.048851969514199
$ws.Range("K7").Value = 1.049322596813124
$ws.Range("L7").Value = 1.053552109625916
$ws.Range("M7").Value = 1.06221361260545
$ws.Range("N7").Value = 1.019636852384269

# Row 8 (bus index 6)
$ws.Range("C8").Value = 1.041036003972119
$ws.Range("D8").Value = 1.043651544782026
$ws.Range("E8").Value = 1.048062736051536
$ws.Range("F8").Value = 1.056457635697644
$ws.Range("J8").Value = 1.046001516165925
$ws.Range("K8").Value = 1.046361850760844
$ws.Range("L8").Value = 1.050760965744564
$ws.Range("M8").Value = 1.059133184507377
$ws.Range("N8").Value = 1.018614591491492

# Row 9 (bus index 7)
$ws.Range("C9").Value = 1.035106927191077
$ws.Range("D9").Value = 1.037921960566552
$ws.Range("E9").Value = 1.042635629575314
$ws.Range("F9").Value = 1.050523089021409
$ws.Range("J9").Value = 1.040911834949643
$ws.Range("K9").Value = 1.041078281382408
$ws.Range("L9").Value = 1.045776711941032
$ws.Range("M9").Value = 1.053639007243573
$ws.Range("N9").Value = 1.01677833938926

# Row 10 (bus index 8)
$ws.Range("C10").Value = 1.031104808891266
$ws.Range("D10").Value = 1.034055502186973
$ws.Range("E10").Value = 1.038971493325862
$ws.Range("F10").Value = 1.046520062579968
$ws.Range("J10").Value = 1.037470649832531
$ws.Range("K10").Value = 1.037508062639325
$ws.Range("L10").Value = 1.042406508592144
$ws.Range("M10").Value = 1.049928490539341
$ws.Range("N10").Value = 1.015529589008377

# Row 11 (bus index 9)
$ws.Range("C11").Value = 1.029359312437538
$ws.Range("D11").Value = 1.032369425976715
$ws.Range("E11").Value = 1.037373225023412
$ws.Range("F11").Value = 1.044774848800068
$ws.Range("J11").Value = 1.035968467098413
$ws.Range("K11").Value = 1.035950041062793
$ws.Range("L11").Value = 1.040935245259233
$ws.Range("M11").Value = 1.048309724469149
$ws.Range("N11").Value = 1.014982781596008

# Row 12 (bus index 10)
$ws.Range("C12").Value = 1.028709000077694
$ws.Range("D12").Value = 1.03174129089536
$ws.Range("E12").Value = 1.036777739503773
$ws.Range("F12").Value = 1.044124745288906
$ws.Range("J12").Value = 1.03540860554721
$ws.Range("K12").Value = 1.035369441849504
$ws.Range("L12").Value = 1.040386897755597
$ws.Range("M12").Value = 1.04770656118122
$ws.Range("N12").Value = 1.014778735749135

# Row 13 (bus index 11)
$ws.Range("C13").Value = 1.028848584021224
$ws.Range("D13").Value = 1.031876112877346
$ws.Range("E13").Value = 1.036905556456136
$ws.Range("F13").Value = 1.044264279696541
$ws.Range("J13").Value = 1.035528784018839
$ws.Range("K13").Value = 1.03549406849484
$ws.Range("L13").Value = 1.040504605077418
$ws.Range("M13").Value = 1.047836027963774
$ws.Range("N13").Value = 1.014822547050643

# Row 14 (bus index 12)
$ws.Range("C14").Value = 1.029305597761362
$ws.Range("D14").Value = 1.032317542296348
$ws.Range("E14").Value = 1.037324039455553
$ws.Range("F14").Value = 1.04472114925531
$ws.Range("J14").Value = 1.035922227565426
$ws.Range("K14").Value = 1.035902087289257
$ws.Range("L14").Value = 1.04088995687591
$ws.Range("M14").Value = 1.048259905569687
$ws.Range("N14").Value = 1.014965934352414

# Row 15 (bus index 13)
$ws.Range("C15").Value = 1.029586917491995
$ws.Range("D15").Value = 1.032589274173725
$ws.Range("E15").Value = 1.037581637873289
$ws.Range("F15").Value = 1.045002393982243
$ws.Range("J15").Value = 1.036164389590152
$ws.Range("K15").Value = 1.03615323000111
$ws.Range("L15").Value = 1.041127137225676
$ws.Range("M15").Value = 1.048520819253089
$ws.Range("N15").Value = 1.015054155135527

# Row 16 (bus index 14)
$ws.Range("C16").Value = 1.031220381217138
$ws.Range("D16").Value = 1.034167145623771
$ws.Range("E16").Value = 1.039077313781741
$ws.Range("F16").Value = 1.046635630581571
$ws.Range("J16").Value = 1.037570083988449
$ws.Range("K16").Value = 1.037611203200614
$ws.Range("L16").Value = 1.042503894736899
$ws.Range("M16").Value = 1.05003566248455
$ws.Range("N16").Value = 1.015565748614914

# Row 17 (bus index 15)
$ws.Range("C17").Value = 1.032241599264275
$ws.Range("D17").Value = 1.035153676581916
$ws.Range("E17").Value = 1.040012341989709
$ws.Range("F17").Value = 1.047656888725778
$ws.Range("J17").Value = 1.038448548867655
$ws.Range("K17").Value = 1.038522469008553
$ws.Range("L17").Value = 1.043364258595997
$ws.Range("M17").Value = 1.050982601642855
$ws.Range("N17").Value = 1.015885011622206

# Row 18 (bus index 16)
$ws.Range("C18").Value = 1.032836052206624
$ws.Range("D18").Value = 1.035727962332573
$ws.Range("E18").Value = 1.040556606213863
$ws.Range("F18").Value = 1.048251430112707
$ws.Range("J18").Value = 1.038959776764373
$ws.Range("K18").Value = 1.039052832485777
$ws.Range("L18").Value = 1.043864946101598
$ws.Range("M18").Value = 1.051533773423067
$ws.Range("N18").Value = 1.016070646239165

# Row 19 (bus index 17)
$ws.Range("C19").Value = 1.033038542831997
$ws.Range("D19").Value = 1.035923587405972
$ws.Range("E19").Value = 1.040741997864389
$ws.Range("F19").Value = 1.048453961861148
$ws.Range("J19").Value = 1.039133896385259
$ws.Range("K19").Value = 1.039233477464123
$ws.Range("L19").Value = 1.044035474663893
$ws.Range("M19").Value = 1.051721513654259
$ws.Range("N19").Value = 1.016133844100685

# Row 20 (bus index 18)
$ws.Range("C20").Value = 1.032132157513226
$ws.Range("D20").Value = 1.035047949654619
$ws.Range("E20").Value = 1.039912138813007
$ws.Range("F20").Value = 1.047547435927113
$ws.Range("J20").Value = 1.038354418929008
$ws.Range("K20").Value = 1.038424819505771
$ws.Range("L20").Value = 1.043272068890828
$ws.Range("M20").Value = 1.050881124672923
$ws.Range("N20").Value = 1.015850818504636

# Row 21 (bus index 19)
$ws.Range("C21").Value = 1.029171073216511
$ws.Range("D21").Value = 1.03218760398818
$ws.Range("E21").Value = 1.037200857333935
$ws.Range("F21").Value = 1.044586664282696
$ws.Range("J21").Value = 1.035806420760359
$ws.Range("K21").Value = 1.035781988357567
$ws.Range("L21").Value = 1.040776532069838
$ws.Range("M21").Value = 1.048135136650864
$ws.Range("N21").Value = 1.01492373641907

# Row 22 (bus index 20)
$ws.Range("C22").Value = 1.027297961370661
$ws.Range("D22").Value = 1.030378445588498
$ws.Range("E22").Value = 1.03548561886827
$ws.Range("F22").Value = 1.042714351047513
$ws.Range("J22").Value = 1.034193461159823
$ws.Range("K22").Value = 1.034109422142072
$ws.Range("L22").Value = 1.039196726786682
$ws.Range("M22").Value = 1.046397705749013
$ws.Range("N22").Value = 1.014335409212803

# Row 23 (bus index 21)
$ws.Range("C23").Value = 1.028292034507777
$ws.Range("D23").Value = 1.031338556130766
$ws.Range("E23").Value = 1.036395920750502
$ws.Range("F23").Value = 1.043707942975484
$ws.Range("J23").Value = 1.035049578893487
$ws.Range("K23").Value = 1.034997137271412
$ws.Range("L23").Value = 1.040035252155885
$ws.Range("M23").Value = 1.047319808003666
$ws.Range("N23").Value = 1.014647815229793

# Row 24 (bus index 22)
$ws.Range("C24").Value = 1.032181613253357
$ws.Range("D24").Value = 1.03509572662457
$ws.Range("E24").Value = 1.039957419778833
$ws.Range("F24").Value = 1.047596896458321
$ws.Range("J24").Value = 1.03839695579567
$ws.Range("K24").Value = 1.038468946701523
$ws.Range("L24").Value = 1.043313728995119
$ws.Range("M24").Value = 1.05092698133427
$ws.Range("N24").Value = 1.015866270712587

# Row 25 (bus index 23)
$ws.Range("C25").Value = 1.03664817334241
$ws.Range("D25").Value = 1.039411178549703
$ws.Range("E25").Value = 1.044046545262018
$ws.Range("F25").Value = 1.052065275295967
$ws.Range("J25").Value = 1.042235876875731
$ws.Range("K25").Value = 1.042452402005795
$ws.Range("L25").Value = 1.047073381547942
$ws.Range("M25").Value = 1.055067549748624
$ws.Range("N25").Value = 1.017257303789771
